$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates for rows 10 and 13-23 (labels/values shifted & two long blocks removed) ---
$ws.Range("B10").Value = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("C10").Value = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("C18").Value = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = @"
LOQ4087 -  Termodinâmica Química Aplicada I  (Requisito fraco)

"@
$ws.Range("C23").Value = @"
LOQ4087 -  Termodinâmica Química Aplicada I  (Requisito fraco)

"@

# --- Row height adjustments to match the new layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30

# --- Remove the trailing row 24 (its content now lives in row 23) ---
$ws.Rows.Item(24).Delete()
